$wb = $excel.ActiveWorkbook

# --- 1. Strategy Status: MarketMaking row (row 6) trade count 1 -> 0 ---
$statusWs = $wb.Worksheets.Item("Strategy Status")
$statusWs.Cells.Item(6, 4).Value = 0

# --- 2. All Trades: close out existing open-trade columns on row 2 and
#        append the new trade as row 3 (history keeps every trade row) ---
$allTradesWs = $wb.Worksheets.Item("All Trades")

# Clear the "open trade" columns (K:Q) on the prior row down to empty text
# cells (quote-prefix blank), matching a closed-out entry.
$allTradesWs.Cells.Item(2, 11).Value = "'"
$allTradesWs.Cells.Item(2, 12).Value = "'"
$allTradesWs.Cells.Item(2, 13).Value = "'"
$allTradesWs.Cells.Item(2, 14).Value = "'"
$allTradesWs.Cells.Item(2, 15).Value = "'"
$allTradesWs.Cells.Item(2, 16).Value = "'"
$allTradesWs.Cells.Item(2, 17).Value = "'"

# Append the new trade as row 3.
$allTradesWs.Cells.Item(3, 1).Value = 2
$allTradesWs.Cells.Item(3, 2).Value = "'2026-02-18"
$allTradesWs.Cells.Item(3, 3).Value = "'10:15:49"
$allTradesWs.Cells.Item(3, 4).Value = "MarketMaking"
$allTradesWs.Cells.Item(3, 5).Value = "DOWN"
$allTradesWs.Cells.Item(3, 6).Value = 0.24
$allTradesWs.Cells.Item(3, 7).Value = "'"
$allTradesWs.Cells.Item(3, 8).Value = "OPEN"
$allTradesWs.Cells.Item(3, 9).Value = 0
$allTradesWs.Cells.Item(3, 10).Value = 0
$allTradesWs.Cells.Item(3, 11).Value = 100
$allTradesWs.Cells.Item(3, 12).Value = 0
$allTradesWs.Cells.Item(3, 13).Value = 0
$allTradesWs.Cells.Item(3, 14).Value = 0.6
$allTradesWs.Cells.Item(3, 15).Value = "Normal spread capture: 202 bps"
$allTradesWs.Cells.Item(3, 16).Value = "'"
$allTradesWs.Cells.Item(3, 17).Value = 0

# --- 3. MarketMaking strategy sheet: overwrite row 2 in place with the
#        latest (now open) trade -- this sheet tracks current state only ---
$mmWs = $wb.Worksheets.Item("MarketMaking")
$mmWs.Cells.Item(2, 1).Value = 2
$mmWs.Cells.Item(2, 3).Value = "'10:15:49"
$mmWs.Cells.Item(2, 5).Value = "DOWN"
$mmWs.Cells.Item(2, 6).Value = 0.24
$mmWs.Cells.Item(2, 7).Value = "'"
$mmWs.Cells.Item(2, 8).Value = "OPEN"
$mmWs.Cells.Item(2, 15).Value = "Normal spread capture: 202 bps"
$mmWs.Cells.Item(2, 16).Value = "'"
$mmWs.Cells.Item(2, 17).Value = 0
